$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 220661.55
$ws.Range("I15").Value = 220661.55
$ws.Range("K15").Value = 661984.6499999999
$ws.Range("M15").Value = -661815.6499999999
$ws.Range("H41").Value = 778.44446
$ws.Range("I41").Value = 773.5714
$ws.Range("J41").Value = 781.5454999999999
$ws.Range("K41").Value = 773.5714
$ws.Range("L41").Value = 781.5454999999999
$ws.Range("M41").Value = -333.5714
$ws.Range("N41").Value = -1661.5455
$ws.Range("H62").Value = 71433570
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 71433570
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 50000
$ws.Range("N65").Value = -56240
$ws.Range("H70").Value = 2080.2
$ws.Range("I70").Value = 1800.6666
$ws.Range("J70").Value = 2499.5
$ws.Range("K70").Value = 5401.9998
$ws.Range("L70").Value = 7498.5
$ws.Range("M70").Value = -5131.9998
$ws.Range("N70").Value = -8038.5
$ws.Range("H73").Value = 2080.2
$ws.Range("I73").Value = 1800.6666
$ws.Range("J73").Value = 2499.5
$ws.Range("K73").Value = 5401.9998
$ws.Range("L73").Value = 7498.5
$ws.Range("M73").Value = -4465.9998
$ws.Range("N73").Value = -9370.5
$ws.Range("H92").Value = 1481.1666
$ws.Range("I92").Value = 2742
$ws.Range("J92").Value = 220.33333
$ws.Range("K92").Value = 2742
$ws.Range("L92").Value = 220.33333
$ws.Range("M92").Value = -1494
$ws.Range("N92").Value = -2716.33333
$ws.Range("H96").Value = 564.73914
$ws.Range("I96").Value = 598.2381
$ws.Range("J96").Value = 213
$ws.Range("K96").Value = 1794.7143
$ws.Range("L96").Value = 639
$ws.Range("M96").Value = -421.7143000000001
$ws.Range("N96").Value = -3385
$ws.Range("H107").Value = 1457.381
$ws.Range("I107").Value = 1652.3572
$ws.Range("J107").Value = 1067.4286
$ws.Range("K107").Value = 1652.3572
$ws.Range("L107").Value = 1067.4286
$ws.Range("M107").Value = 267.6428000000001
$ws.Range("N107").Value = -4907.4286
$ws.Range("H125").Value = 3086.1875
$ws.Range("I125").Value = 4304.857
$ws.Range("J125").Value = 2138.3333
$ws.Range("K125").Value = 38743.713
$ws.Range("L125").Value = 19244.9997
$ws.Range("M125").Value = -36283.713
$ws.Range("N125").Value = -24164.9997
$ws.Range("H137").Value = 13979.056
$ws.Range("I137").Value = 12142.714
$ws.Range("K137").Value = 36428.142
$ws.Range("M137").Value = -33878.142
$ws.Range("H138").Value = 4893.8286
$ws.Range("J138").Value = 6166.8887
$ws.Range("L138").Value = 18500.6661
$ws.Range("N138").Value = -28780.6661
$ws.Range("H141").Value = 2097.8333
$ws.Range("I141").Value = 2094.25
$ws.Range("K141").Value = 6282.75
$ws.Range("M141").Value = -1102.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 8353916.5
$ws.Range("I10").Value = 25017500
$ws.Range("J10").Value = 22125
$ws.Range("K10").Value = 25017500
$ws.Range("L10").Value = 22125
$ws.Range("M10").Value = -25017330
$ws.Range("N10").Value = -22465
$ws.Range("H11").Value = 50000000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H32").Value = 1324633.5
$ws.Range("I32").Value = 1085.2972
$ws.Range("J32").Value = 11118890
$ws.Range("K32").Value = 1085.2972
$ws.Range("L32").Value = 11118890
$ws.Range("M32").Value = -798.2972
$ws.Range("N32").Value = -11119464
$ws.Range("H45").Value = 2859.8
$ws.Range("I45").Value = 2859.8
$ws.Range("K45").Value = 2859.8
$ws.Range("M45").Value = -2482.8
$ws.Range("H61").Value = 4822.6924
$ws.Range("I61").Value = 5346.5293
$ws.Range("J61").Value = 3833.2222
$ws.Range("K61").Value = 5346.5293
$ws.Range("L61").Value = 3833.2222
$ws.Range("M61").Value = -5134.5293
$ws.Range("N61").Value = -4257.2222
$ws.Range("H74").Value = 2718.9143
$ws.Range("I74").Value = 2874.037
$ws.Range("K74").Value = 2874.037
$ws.Range("M74").Value = -2000.037
$ws.Range("H77").Value = 2718.9143
$ws.Range("I77").Value = 2874.037
$ws.Range("K77").Value = 14370.185
$ws.Range("M77").Value = -10002.185
$ws.Range("H97").Value = 745.3333
$ws.Range("I97").Value = 548
$ws.Range("K97").Value = 548
$ws.Range("M97").Value = -52
$ws.Range("H122").Value = 5999.8335
$ws.Range("I122").Value = 5499.75
$ws.Range("K122").Value = 16499.25
$ws.Range("M122").Value = -14049.25
$ws.Range("H132").Value = 923878.0600000001
$ws.Range("I132").Value = 1168162
$ws.Range("J132").Value = 129955.5
$ws.Range("K132").Value = 3504486
$ws.Range("L132").Value = 389866.5
$ws.Range("M132").Value = -3501956
$ws.Range("N132").Value = -394926.5
$ws.Range("H136").Value = 4822.6924
$ws.Range("I136").Value = 5346.5293
$ws.Range("J136").Value = 3833.2222
$ws.Range("K136").Value = 16039.5879
$ws.Range("L136").Value = 11499.6666
$ws.Range("M136").Value = -13489.5879
$ws.Range("N136").Value = -16599.6666

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 6394.8857
$ws.Range("I22").Value = 2844
$ws.Range("J22").Value = 27700.2
$ws.Range("K22").Value = 2844
$ws.Range("L22").Value = 27700.2
$ws.Range("M22").Value = -2671
$ws.Range("N22").Value = -28046.2
$ws.Range("H94").Value = 30795.27
$ws.Range("I94").Value = 596.7778
$ws.Range("J94").Value = 59404.367
$ws.Range("K94").Value = 596.7778
$ws.Range("L94").Value = 59404.367
$ws.Range("M94").Value = -145.7778
$ws.Range("N94").Value = -60306.367
$ws.Range("H134").Value = 1255311.8
$ws.Range("I134").Value = 1474219.8
$ws.Range("K134").Value = 4422659.4
$ws.Range("M134").Value = -4420124.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5061.5864
$ws.Range("I31").Value = 2176.3
$ws.Range("J31").Value = 6580.1577
$ws.Range("K31").Value = 2176.3
$ws.Range("L31").Value = 6580.1577
$ws.Range("M31").Value = -1881.3
$ws.Range("N31").Value = -7170.1577
$ws.Range("H34").Value = 5061.5864
$ws.Range("I34").Value = 2176.3
$ws.Range("J34").Value = 6580.1577
$ws.Range("K34").Value = 2176.3
$ws.Range("L34").Value = 6580.1577
$ws.Range("M34").Value = -1974.3
$ws.Range("N34").Value = -6984.1577
$ws.Range("H58").Value = 18523210
$ws.Range("I58").Value = 25003166
$ws.Range("K58").Value = 25003166
$ws.Range("M58").Value = -25002963
$ws.Range("H105").Value = 200002530
$ws.Range("H132").Value = 4951.7905
$ws.Range("I132").Value = 4463.355
$ws.Range("K132").Value = 13390.065
$ws.Range("M132").Value = -10860.065
$ws.Range("H136").Value = 18523210
$ws.Range("I136").Value = 25003166
$ws.Range("K136").Value = 75009498
$ws.Range("M136").Value = -75006948

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15916.5
$ws.Range("I70").Value = 8222
$ws.Range("K70").Value = 8222
$ws.Range("M70").Value = -7952
$ws.Range("H73").Value = 15916.5
$ws.Range("I73").Value = 8222
$ws.Range("K73").Value = 8222
$ws.Range("M73").Value = -7286
$ws.Range("H97").Value = 36471.145
$ws.Range("I97").Value = 812.85
$ws.Range("K97").Value = 812.85
$ws.Range("M97").Value = -316.85
$ws.Range("H107").Value = 699.61536
$ws.Range("J107").Value = 757.6
$ws.Range("L107").Value = 757.6
$ws.Range("N107").Value = -4597.6
$ws.Range("H126").Value = 17246234
$ws.Range("I126").Value = 29414930
$ws.Range("J126").Value = 7246.5
$ws.Range("K126").Value = 88244790
$ws.Range("L126").Value = 21739.5
$ws.Range("M126").Value = -88242320
$ws.Range("N126").Value = -26679.5
$ws.Range("H132").Value = 16132550
$ws.Range("J132").Value = 4036.5715
$ws.Range("L132").Value = 12109.7145
$ws.Range("N132").Value = -17169.7145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5956.5713
$ws.Range("I40").Value = 5062.25
$ws.Range("K40").Value = 5062.25
$ws.Range("M40").Value = -4926.25
$ws.Range("H55").Value = 3007.12
$ws.Range("I55").Value = 698.2857
$ws.Range("J55").Value = 5945.636
$ws.Range("K55").Value = 698.2857
$ws.Range("L55").Value = 5945.636
$ws.Range("M55").Value = -525.2857
$ws.Range("N55").Value = -6291.636
$ws.Range("H122").Value = 29950
$ws.Range("I122").Value = 38933.332
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 116799.996
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -114349.996
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 11279.25
$ws.Range("I132").Value = 9752.5
$ws.Range("K132").Value = 29257.5
$ws.Range("M132").Value = -26727.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 8340968
$ws.Range("I17").Value = 8340968
$ws.Range("K17").Value = 8340968
$ws.Range("M17").Value = -8340796
$ws.Range("H132").Value = 5747.3438
$ws.Range("J132").Value = 9051.429
$ws.Range("L132").Value = 27154.287
$ws.Range("N132").Value = -32214.287

Write-Host "Edit complete"